$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update initial condition values
$ws.Range("I2").Value = 0.0
$ws.Range("K2").Value = -5.0
$ws.Range("L2").Value = 0.0

# Row 3: keep Timestep label (A3) and B3 (0.0), clear out all the simulated
# trajectory columns so the row only carries the initial condition marker.
$ws.Range("D3:G3").ClearContents()
$ws.Range("I3:L3").ClearContents()
$ws.Range("N3:O3").ClearContents()

# Row 4: same cleanup, plus reset B4 back to 0.0
$ws.Range("B4").Value = 0.0
$ws.Range("D4:G4").ClearContents()
$ws.Range("I4:L4").ClearContents()
$ws.Range("N4:O4").ClearContents()
